# Update "paises.xlsx" (sheet "Pais") with the refreshed COVID-19 snapshot:
#  - bump the "last updated" note in A1 to the new capture time
#  - a handful of countries changed rank (their total-case count crossed a
#    neighbour's), so the row that used to show one country's data now shows
#    another's - update column A for those rows accordingly
#  - write the refreshed Casos totales / Nuevos casos / Casos activos /
#    Recuperados / Casos criticos / Muertes hoy / Muertes figures (B:H) for
#    every row whose numbers moved
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp note (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 09:47"

# --- Country-name cell swaps caused by re-sorting the source data by case count ---
$ws.Range("A60").Value = "Armenia"
$ws.Range("A61").Value = "Uzbekistan"
$ws.Range("A73").Value = "Hungria"
$ws.Range("A74").Value = "Azerbaiyan"
$ws.Range("A75").Value = "Kenia"
$ws.Range("A96").Value = "Georgia"
$ws.Range("A97").Value = "Noruega"
$ws.Range("A98").Value = "Zambia"
$ws.Range("A99").Value = "Senegal"
$ws.Range("A141").Value = "Estonia"
$ws.Range("A142").Value = "Mayotte"
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("A217").Value = "Montserrat"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B7").Value = 1384235
$ws.Range("C7").Value = 14922
$ws.Range("D7").Value = 1065199
$ws.Range("E7").Value = 295034
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 279
$ws.Range("H7").Value = 24002

$ws.Range("B60").Value = 63000
$ws.Range("C60").Value = 1540
$ws.Range("D60").Value = 47925
$ws.Range("E60").Value = 14008
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 11
$ws.Range("H60").Value = 1067

$ws.Range("B61").Value = 62684
$ws.Range("C61").Value = 96
$ws.Range("D61").Value = 59756
$ws.Range("E61").Value = 2408
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 520

$ws.Range("B73").Value = 44816
$ws.Range("C73").Value = 1791
$ws.Range("D73").Value = 13580
$ws.Range("E73").Value = 30127
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 24
$ws.Range("H73").Value = 1109

$ws.Range("B74").Value = 43789
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 39800
$ws.Range("E74").Value = 3368
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 621

$ws.Range("B75").Value = 43580
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 31648
$ws.Range("E75").Value = 11119
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 813

$ws.Range("B76").Value = 40141
$ws.Range("C76").Value = 68
$ws.Range("D76").Value = 33561
$ws.Range("E76").Value = 5092
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 1488

$ws.Range("B96").Value = 16285
$ws.Range("C96").Value = 958
$ws.Range("D96").Value = 7827
$ws.Range("E96").Value = 8330
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 4
$ws.Range("H96").Value = 128

$ws.Range("B97").Value = 16272
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 11863
$ws.Range("E97").Value = 4131
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 278

$ws.Range("B98").Value = 15659
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 14899
$ws.Range("E98").Value = 414
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 346

$ws.Range("B99").Value = 15368
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 13704
$ws.Range("E99").Value = 1347
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 317

$ws.Range("B141").Value = 4052
$ws.Range("C141").Value = 35
$ws.Range("D141").Value = 3198
$ws.Range("E141").Value = 786
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 68

$ws.Range("B142").Value = 4030
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 2964
$ws.Range("E142").Value = 1023
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 43

$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
